$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the table with the 2020 column (K), mirroring the formatting of the
# existing 2019 column (J) and filling in the new figures.
$ws.Range("J2:J8").Copy()
$ws.Range("K2:K8").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("K3").Value = 2020
$ws.Range("K4").Value = 0
$ws.Range("K5").Value = 48.2
$ws.Range("K6").Value = 19.3
$ws.Range("K7").Value = 24.2
$ws.Range("K8").Value = 8.3000000000000007

# The header row's wrapped title now renders a bit shorter once recalculated.
$ws.Rows.Item(1).RowHeight = 63.75

# Restore the cursor position left by the author after editing.
$ws.Range("J22").Select()
